# This script applies the diff described:
#  - Row 5 and Row 6 swap their full species-observation data
#    (row 5 becomes the old row 6 content, row 6 becomes the old row 5 content)
#  - Row 7 and Row 8 swap only their Id (A), Ost (Q), Nord (R) coordinates,
#    and the "Publik kommentar" (AC) comment moves from row 8 to row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 5 <-> Row 6 full swap ----

# New row 5 (previously row 6's content: Motaggsvamp / Sarcodon squamosus)
$ws.Range("A5").Value = 111934086
$ws.Range("B5").Value = 90689
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 5966
$ws.Range("F5").Value = "Motaggsvamp"
$ws.Range("G5").Value = "Sarcodon squamosus"
$ws.Range("H5").Value = "(Schaeff.) Quél."

# I5 needs to hold the text "1" (stored as text, not a number) -- force with a
# leading apostrophe, then reset the style so no stray quote-prefix format
# lingers on the cell.
$ws.Range("I5").Value = "'1"
$ws.Range("I5").Style = "Normal"

$ws.Range("J5").Value = "fruktkroppar"

# K5 becomes an empty (but still text-typed) cell.
$ws.Range("K5").Value = "'"
$ws.Range("K5").Style = "Normal"

# L5 had content before (blank placeholder) and should end up fully empty.
$ws.Range("L5").Value = ""

$ws.Range("P5").Value = "Tallskogen N om Dye, I2-Skogen, Vrm"
$ws.Range("Q5").Value = 413681.2082122188
$ws.Range("R5").Value = 6586805.223123537

# New row 6 (previously row 5's content: Flagellkvastmossa / Dicranum flagellare)
$ws.Range("A6").Value = 111934050
$ws.Range("B6").Value = 93289
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 2170
$ws.Range("F6").Value = "Flagellkvastmossa"
$ws.Range("G6").Value = "Dicranum flagellare"
$ws.Range("H6").Value = "Hedw."

# I6 and J6 become empty (but still text-typed) cells.
$ws.Range("I6").Value = "'"
$ws.Range("I6").Style = "Normal"
$ws.Range("J6").Value = "'"
$ws.Range("J6").Style = "Normal"

$ws.Range("K6").Value = "med groddkorn"

# L6 newly appears as an empty (but still text-typed) cell.
$ws.Range("L6").Value = "'"
$ws.Range("L6").Style = "Normal"

$ws.Range("P6").Value = "Skogen N om Dye, I2-Skogen, Vrm"
$ws.Range("Q6").Value = 413637.9321653559
$ws.Range("R6").Value = 6587076.603947581

# ---- Row 7 <-> Row 8 partial swap (Id, Ost, Nord, Publik kommentar) ----

$ws.Range("A7").Value = 111934059
$ws.Range("Q7").Value = 413639.6308819132
$ws.Range("R7").Value = 6586793.951973591
$ws.Range("AC7").Value = "Rätt riklig längs stigen"

$ws.Range("A8").Value = 111934066
$ws.Range("Q8").Value = 413590.3038565172
$ws.Range("R8").Value = 6586912.201658082
$ws.Range("AC8").Value = ""
